# The deck had two slides near the end that needed to be removed:
#   slide 25 "Sub-Title 2"  (title placeholder + picture + table + captions)
#   slide 26 "Reference"    (title + "Use latest APA style referencing" + URL)
# After removing them, the former slide 27 "Thank you" becomes the new
# (last) slide 25, and the presentation's notes-master / slide-id
# relationship bookkeeping shifts down accordingly.
#
# Delete from the higher index first so the lower index ("Sub-Title 2")
# still refers to the same slide after the first deletion.

$p = $ppt.ActivePresentation

$p.Slides.Item(26).Delete()   # "Reference"
$p.Slides.Item(25).Delete()   # "Sub-Title 2"
